$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "3 V 0.3" : break out row for LXCHEM / Laxmi Organic Industries Ltd
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("3 V 0.3")

# E10 was stored as text "543277" - convert it to a real number
$ws1.Cells.Item(10, 5).Value = 543277

# Append new row 11 with the newly split-out data point
$ws1.Cells.Item(11, 1).Value = "20/06/2024 09:45:37"
$ws1.Cells.Item(11, 2).Value = 1
$ws1.Cells.Item(11, 3).Value = "LXCHEM"
$ws1.Cells.Item(11, 4).Value = "Laxmi Organic Industries Ltd"
$ws1.Cells.Item(11, 5).NumberFormat = "@"
$ws1.Cells.Item(11, 5).Value = "543277"
$ws1.Cells.Item(11, 6).Value = 3.56
$ws1.Cells.Item(11, 7).Value = 263.5
$ws1.Cells.Item(11, 8).Value = 6238938

# ---------------------------------------------------------------------------
# Sheet "DND 3 V 0.3" : break out row for IBREALEST / Indiabulls Real Estate
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("DND 3 V 0.3")

# E6 was stored as text "532832" - convert it to a real number
$ws2.Cells.Item(6, 5).Value = 532832

# Append new row 7 with the newly split-out data point
$ws2.Cells.Item(7, 1).Value = "20/06/2024 09:45:37"
$ws2.Cells.Item(7, 2).Value = 1
$ws2.Cells.Item(7, 3).Value = "IBREALEST"
$ws2.Cells.Item(7, 4).Value = "Indiabulls Real Estate Limited"
$ws2.Cells.Item(7, 5).NumberFormat = "@"
$ws2.Cells.Item(7, 5).Value = "532832"
$ws2.Cells.Item(7, 6).Value = 12.9
$ws2.Cells.Item(7, 7).Value = 154.58
$ws2.Cells.Item(7, 8).Value = 84581155
